$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title replacements (heading + bold recap at the end)
Replace-Text "Play Aztec Magic Deluxe for Free - Review and Features" "Play Aztec Magic Deluxe for Free"

# "What we like" bullet list
Replace-Text "Beautifully designed with Aztec-themed symbols" "Beautifully designed with attention to detail"
Replace-Text "Variable pay lines and Autoplay for up to 1,000 spins" "Varied gameplay features and adjustable betting options"
Replace-Text "Aztec warrior symbol triples winnings and pays up to 5,000 times the bet" "Themed symbols pay homage to Aztec culture"
Replace-Text "Good RTP value of 96.96%" "Medium volatility and high RTP for balanced risk and reward"

# "What we don't like" bullet list
Replace-Text "Only one special symbol - the Wild - which also serves as the highest paying symbol" "Limited number of free spins triggered by Scatter symbol"
Replace-Text "Maximum cost per spin is only €1" "No progressive jackpot feature"

# Final italic summary paragraph
Replace-Text "Discover the Aztec-themed slot game Aztec Magic Deluxe and play it for free. Read our review of the features and gameplay of this beautifully designed game." "Read an unbiased review of Aztec Magic Deluxe and play it for free. Experience the beauty of Aztec culture in this captivating slot game."
